$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header cells in columns C and D (no longer used).
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = ""

$colA = @(
    'everton mendozaholmes',
    '0FD8AD9C60',
    '2017/01/19 15:03',
    '2017/01/19 15:03',
    '2017/01/19 15:04',
    '2017/01/19 15:04',
    '2017/01/19 15:05',
    '2017/01/19 15:06',
    '2017/01/19 15:06',
    '2017/01/19 15:06',
    '2017/01/19 15:06',
    '2017/01/19 15:15',
    '2017/01/19 15:16',
    '2017/01/19 15:19',
    '2017/01/19 15:24',
    '2017/01/19 15:25',
    '2017/01/19 15:25',
    '2017/01/19 15:26',
    '2017/01/19 15:27',
    '2017/01/19 15:27',
    '2017/01/19 15:27',
    '2017/01/19 15:29',
    '2017/01/19 15:30',
    '2017/01/19 15:30',
    '2017/01/19 15:30',
    '2017/01/19 15:30',
    '2017/01/19 15:30',
    '2017/01/19 15:31',
    '2017/01/19 15:31',
    '2017/01/19 15:31',
    '2017/01/19 15:31',
    '2017/01/19 15:32',
    '2017/01/19 15:32',
    '2017/01/19 15:32',
    '2017/01/19 15:34',
    '2017/01/19 15:34',
    '2017/01/19 15:43',
    '2017/01/19 15:44',
    '2017/01/19 15:45',
    '2017/01/19 15:46',
    '2017/01/19 15:50',
    '2017/01/19 15:51',
    '2017/01/19 15:52',
    '2017/01/19 15:53',
    '2017/01/19 15:53',
    '2017/01/19 15:54',
    '2017/01/19 15:54',
    '2017/01/19 15:56',
    '2017/01/19 15:56',
    '2017/01/19 15:56',
    '2017/01/19 15:57',
    '2017/01/19 16:01',
    '2017/01/19 16:01',
    '2017/01/19 16:04',
    '2017/01/19 16:04',
    '2017/01/19 17:28',
    '2017/01/19 17:35',
    '2017/01/19 17:44',
    '2017/01/19 17:44',
    '2017/01/19 17:44',
    '2017/01/19 17:46',
    '2017/01/19 17:46',
    '2017/01/19 17:46',
    '2017/01/19 18:00',
    '2017/01/19 18:00',
    '2017/01/19 18:00',
    '2017/01/19 18:09',
    '2017/01/19 18:09'
)
$colB = @(
    'Gary Tsai',
    'CFD893A460',
    '2017/01/19 18:59',
    '2017/01/19 19:00',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    '',
    ''
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    if ($colB[$i] -ne '') {
        $ws.Cells.Item($r, 2).Value = $colB[$i]
    }
}
